# MOS-23045: Update Master Data as per 22 April Changes
# Adds new POA/RNC rows for apptyp_code 3,4,7,8,11,12,15 and refreshes the
# sheet's AutoFilter range over the full data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-apply AutoFilter across the existing table BEFORE appending rows so the
#     filter range stays anchored to the original A1:G57 extent. ---
[void]$ws.Range("A1:G57").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='master-valid_document'!`$A`$1:`$G`$57")
$filterName.Visible = $false

# --- New data rows (apptyp_code, doccat_code, doctyp_code, lang_code, is_active, cr_by, cr_dtimes) ---
$newRows = @(
    @(58, 3),
    @(59, 4),
    @(60, 7),
    @(61, 8),
    @(62, 11),
    @(63, 12),
    @(64, 15)
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    $apptyp = $entry[1]

    $ws.Range("A$r").Value = $apptyp
    $ws.Range("A$r").NumberFormat = "000"

    $ws.Range("B$r").Value = "POA"
    $ws.Range("C$r").Value = "RNC"
    $ws.Range("D$r").Value = "eng"
    $ws.Range("E$r").Value = $true
    $ws.Range("F$r").Value = "superadmin"
    $ws.Range("G$r").Value = "now()"
}

# --- Header row formatting touch-up (matches the numeric-style cleanup made alongside the new rows) ---
$ws.Range("A1").NumberFormat = "000"

# --- Selection left where the editor's cursor ended up ---
[void]$ws.Range("H2").Select()
